# Estadisticos Segundo Parcial Sin Ameca
# Adds a new "Ingles II - 2APV" group and a new "Ingles IV - 4AEV" group to
# each of the three statistics sheets, and adds the corresponding makeup
# ("Rescatables") students, re-sorted by number of failed subjects.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets "Estadisticos 1P", "Estadisticos 2P" y "Estadisticos Final"
# ---------------------------------------------------------------------

$statSheets = @(
  @{ Name = "Estadisticos 1P";
     Rows = @(
       @("Ingles II", "2APV",  23, 0,  7, 16, 69.56999999999999, 6.5),
       @("Ingles IV", "4AEV",  19, 0, 12,  7, 36.84,              5.5),
       @("Ingles IV", "4ALCV", 25, 0, 12, 13, 52,                 6.3),
       @("Ingles IV", "4APV",  13, 0,  5,  8, 61.54,              6.2),
       @("Ingles IV", "4ARHV", 14, 0,  0, 14, 100,                9.300000000000001),
       @("Ingles IV", "4ASV",  10, 0,  2,  8, 80,                 6.8)
     ) },
  @{ Name = "Estadisticos 2P";
     Rows = @(
       @("Ingles II", "2APV",  23, 0,  5, 18, 78.26000000000001, 6.5),
       @("Ingles IV", "4AEV",  19, 0,  9, 10, 52.63,              5.5),
       @("Ingles IV", "4ALCV", 25, 0,  7, 18, 72,                 6.3),
       @("Ingles IV", "4APV",  13, 0,  5,  8, 61.54,              6.2),
       @("Ingles IV", "4ARHV", 14, 0,  0, 14, 100,                9.300000000000001),
       @("Ingles IV", "4ASV",  10, 0,  1,  9, 90,                 6.8)
     ) },
  @{ Name = "Estadisticos Final";
     Rows = @(
       @("Ingles II", "2APV",  23, 0,  5, 18, 78.26000000000001, 7),
       @("Ingles IV", "4AEV",  19, 0,  9, 10, 52.63,              6),
       @("Ingles IV", "4ALCV", 25, 0,  7, 18, 72,                 7.1),
       @("Ingles IV", "4APV",  13, 0,  5,  8, 61.54,              6.6),
       @("Ingles IV", "4ARHV", 14, 0,  0, 14, 100,                9.1),
       @("Ingles IV", "4ASV",  10, 0,  1,  9, 90,                 7.9)
     ) }
)

foreach ($entry in $statSheets) {
  $ws = $wb.Worksheets.Item($entry.Name)
  $r = 2
  foreach ($row in $entry.Rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
  }
}

# ---------------------------------------------------------------------
# Sheet "Rescatables" - list of students who must retake the exam,
# sorted descending by number of failed subjects (column G)
# ---------------------------------------------------------------------

$rescatables = @(
  @(24330051920393, "MUÑOZ",      "REYES",      "ERWIN ISRAEL",      "Ingles II", "2APV",  4),
  @(23330051920211, "VAZQUEZ",    "CARRILLO",   "DIEGO ARMANDO",     "Ingles IV", "4AEV",  4),
  @(23330051920301, "ORTIZ",      "CORTES",     "CARLOS",            "Ingles IV", "4APV",  4),
  @(22330051920389, "FLORES",     "LAGUNA",     "JOSE ANTONIO",      "Ingles IV", "4AEV",  3),
  @(23330051920332, "RODRIGUEZ",  "SUAREZ",     "SERGIO JOSUE",      "Ingles IV", "4AEV",  3),
  @(23330051920238, "HUERTA",     "ALCANTARA",  "JHON STEVE",        "Ingles IV", "4ALCV", 3),
  @(23330051920264, "SAN JUAN",   "AGUILAR",    "KARINA MONSERRATH", "Ingles IV", "4ALCV", 3),
  @(23330051920296, "LUCIANO",    "CAMPOS",     "KEVIN SANTIAGO",    "Ingles IV", "4APV",  3),
  @(23330051920329, "JIMENEZ",    "CIRUELO",    "SABDY",             "Ingles IV", "4AEV",  2),
  @(23330051920324, "JUAREZ",     "LIBRADO",    "ARMANDO GABRIEL",   "Ingles IV", "4AEV",  2),
  @(23330051920203, "PEREZ",      "DE JESUS",   "LUIS FABIAN",       "Ingles IV", "4AEV",  2),
  @(23330051920220, "BARRAGAN",   "MACUIXTLE",  "MARIA FERNANDA",    "Ingles IV", "4ALCV", 2),
  @(23330051920228, "GARCIA",     "ROSALES",    "REGINA DAYTRI",     "Ingles IV", "4ALCV", 2),
  @(23330051920267, "TEXCAHUA",   "DE LA CRUZ", "KARINA",            "Ingles IV", "4ALCV", 2),
  @(23330051920309, "PELAYO",     "TORRES",     "IVAN JESUS",        "Ingles IV", "4APV",  2),
  @(23330051920290, "FLORES",     "TINOCO",     "JULIO ALBERTO",     "Ingles IV", "4ASV",  2),
  @(22330051920240, "GARCIA",     "MARQUEZ",    "AMALIA PAULINA",    "Ingles IV", "4ALCV", 1),
  @(23330051920271, "IXMATLAHUA", "HERNANDEZ",  "FERNANDA YAMILET",  "Ingles IV", "4ALCV", 1),
  @(23330051920294, "FLORES",     "SANCHEZ",    "LUIS PABLO",        "Ingles IV", "4APV",  1),
  @(23330051920363, "ROMERO",     "MARTINEZ",   "AARON",             "Ingles IV", "4APV",  1)
)

$wsR = $wb.Worksheets.Item("Rescatables")
$r = 2
foreach ($row in $rescatables) {
  $wsR.Cells.Item($r, 1).Value = $row[0]
  $wsR.Cells.Item($r, 2).Value = $row[1]
  $wsR.Cells.Item($r, 3).Value = $row[2]
  $wsR.Cells.Item($r, 4).Value = $row[3]
  $wsR.Cells.Item($r, 5).Value = $row[4]
  $wsR.Cells.Item($r, 6).Value = $row[5]
  $wsR.Cells.Item($r, 7).Value = $row[6]
  $r = $r + 1
}
